# Auto-generated edit script: applies market-price / leve-profit updates
# from the scheduled Kujata data refresh across all 8 crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 580.4286
$ws.Cells.Item(2, 9).Value = 580.4286
$ws.Cells.Item(2, 11).Value = 580.4286
$ws.Cells.Item(2, 13).Value = -467.4286
$ws.Cells.Item(15, 8).Value = 883.1
$ws.Cells.Item(15, 9).Value = 883.1
$ws.Cells.Item(15, 11).Value = 2649.3
$ws.Cells.Item(15, 13).Value = -2480.3
$ws.Cells.Item(29, 8).Value = 1476.3334
$ws.Cells.Item(29, 10).Value = 1803.75
$ws.Cells.Item(29, 12).Value = 5411.25
$ws.Cells.Item(29, 14).Value = -5973.25
$ws.Cells.Item(43, 8).Value = 7939951
$ws.Cells.Item(43, 10).Value = 18525186
$ws.Cells.Item(43, 12).Value = 18525186
$ws.Cells.Item(43, 14).Value = -18525324
$ws.Cells.Item(51, 8).Value = 909.0909
$ws.Cells.Item(51, 9).Value = 666.6667
$ws.Cells.Item(51, 10).Value = 1000
$ws.Cells.Item(51, 11).Value = 666.6667
$ws.Cells.Item(51, 12).Value = 1000
$ws.Cells.Item(51, 13).Value = -182.6667
$ws.Cells.Item(51, 14).Value = -1968
$ws.Cells.Item(127, 8).Value = 2673.7
$ws.Cells.Item(127, 9).Value = 887
$ws.Cells.Item(127, 10).Value = 2872.2222
$ws.Cells.Item(127, 11).Value = 2661
$ws.Cells.Item(127, 12).Value = 8616.6666
$ws.Cells.Item(127, 13).Value = 2299
$ws.Cells.Item(127, 14).Value = -18536.6666
$ws.Cells.Item(128, 8).Value = 25000
$ws.Cells.Item(128, 10).Value = 25000
$ws.Cells.Item(128, 12).Value = 25000
$ws.Cells.Item(128, 14).Value = -34960
$ws.Cells.Item(129, 8).Value = 918.375
$ws.Cells.Item(129, 9).Value = 797.6667
$ws.Cells.Item(129, 10).Value = 990.8
$ws.Cells.Item(129, 11).Value = 2393.0001
$ws.Cells.Item(129, 12).Value = 2972.4
$ws.Cells.Item(129, 13).Value = 2606.9999
$ws.Cells.Item(129, 14).Value = -12972.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(9, 8).Value = 10000
$ws.Cells.Item(9, 10).Value = 10000
$ws.Cells.Item(9, 12).Value = 10000
$ws.Cells.Item(9, 14).Value = -10340
$ws.Cells.Item(20, 8).Value = 10000
$ws.Cells.Item(20, 10).Value = 10000
$ws.Cells.Item(20, 12).Value = 10000
$ws.Cells.Item(20, 14).Value = -10540
$ws.Cells.Item(32, 8).Value = 21475.547
$ws.Cells.Item(32, 9).Value = 10378.164
$ws.Cells.Item(32, 10).Value = 36007.832
$ws.Cells.Item(32, 11).Value = 10378.164
$ws.Cells.Item(32, 12).Value = 36007.832
$ws.Cells.Item(32, 13).Value = -10091.164
$ws.Cells.Item(32, 14).Value = -36581.832
$ws.Cells.Item(37, 8).Value = 21629.334
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(61, 8).Value = 43479790
$ws.Cells.Item(61, 9).Value = 52632840
$ws.Cells.Item(61, 11).Value = 52632840
$ws.Cells.Item(61, 13).Value = -52632628
$ws.Cells.Item(102, 8).Value = 8336131
$ws.Cells.Item(102, 9).Value = 8336131
$ws.Cells.Item(102, 11).Value = 8336131
$ws.Cells.Item(102, 13).Value = -8334509
$ws.Cells.Item(136, 8).Value = 43479790
$ws.Cells.Item(136, 9).Value = 52632840
$ws.Cells.Item(136, 11).Value = 157898520
$ws.Cells.Item(136, 13).Value = -157895970
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 38462924
$ws.Cells.Item(105, 10).Value = 2750
$ws.Cells.Item(105, 12).Value = 2750
$ws.Cells.Item(105, 14).Value = -6244
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 307.5
$ws.Cells.Item(7, 9).Value = 272
$ws.Cells.Item(7, 10).Value = 366.66666
$ws.Cells.Item(7, 11).Value = 272
$ws.Cells.Item(7, 12).Value = 366.66666
$ws.Cells.Item(7, 13).Value = -159
$ws.Cells.Item(7, 14).Value = -592.66666
$ws.Cells.Item(31, 8).Value = 1562.9678
$ws.Cells.Item(31, 9).Value = 1515.614
$ws.Cells.Item(31, 11).Value = 1515.614
$ws.Cells.Item(31, 13).Value = -1220.614
$ws.Cells.Item(34, 8).Value = 1562.9678
$ws.Cells.Item(34, 9).Value = 1515.614
$ws.Cells.Item(34, 11).Value = 1515.614
$ws.Cells.Item(34, 13).Value = -1313.614
$ws.Cells.Item(105, 8).Value = 735.6667
$ws.Cells.Item(105, 9).Value = 640.63635
$ws.Cells.Item(105, 11).Value = 640.63635
$ws.Cells.Item(105, 13).Value = 1106.36365
$ws.Cells.Item(122, 8).Value = 1922
$ws.Cells.Item(122, 9).Value = 1897.6666
$ws.Cells.Item(122, 11).Value = 5692.9998
$ws.Cells.Item(122, 13).Value = -3242.9998
$ws.Cells.Item(132, 8).Value = 1549.0209
$ws.Cells.Item(132, 9).Value = 1248.4286
$ws.Cells.Item(132, 10).Value = 2358.3076
$ws.Cells.Item(132, 11).Value = 3745.2858
$ws.Cells.Item(132, 12).Value = 7074.9228
$ws.Cells.Item(132, 13).Value = -1215.2858
$ws.Cells.Item(132, 14).Value = -12134.9228
$ws.Cells.Item(134, 8).Value = 9260543
$ws.Cells.Item(134, 9).Value = 1197.0278
$ws.Cells.Item(134, 10).Value = 27779234
$ws.Cells.Item(134, 11).Value = 3591.0834
$ws.Cells.Item(134, 12).Value = 83337702
$ws.Cells.Item(134, 13).Value = -1056.0834
$ws.Cells.Item(134, 14).Value = -83342772
$ws.Cells.Item(141, 8).Value = 284950.6
$ws.Cells.Item(141, 10).Value = 284950.6
$ws.Cells.Item(141, 12).Value = 284950.6
$ws.Cells.Item(141, 14).Value = -295310.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 23846162
$ws.Cells.Item(131, 10).Value = 46527
$ws.Cells.Item(131, 12).Value = 139581
$ws.Cells.Item(131, 14).Value = -149661
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 19568908
$ws.Cells.Item(70, 9).Value = 17860892
$ws.Cells.Item(70, 10).Value = 22225824
$ws.Cells.Item(70, 11).Value = 17860892
$ws.Cells.Item(70, 12).Value = 22225824
$ws.Cells.Item(70, 13).Value = -17860622
$ws.Cells.Item(70, 14).Value = -22226364
$ws.Cells.Item(73, 8).Value = 19568908
$ws.Cells.Item(73, 9).Value = 17860892
$ws.Cells.Item(73, 10).Value = 22225824
$ws.Cells.Item(73, 11).Value = 17860892
$ws.Cells.Item(73, 12).Value = 22225824
$ws.Cells.Item(73, 13).Value = -17859956
$ws.Cells.Item(73, 14).Value = -22227696
$ws.Cells.Item(102, 8).Value = 2774.6453
$ws.Cells.Item(102, 9).Value = 1797.95
$ws.Cells.Item(102, 10).Value = 4550.4546
$ws.Cells.Item(102, 11).Value = 1797.95
$ws.Cells.Item(102, 12).Value = 4550.4546
$ws.Cells.Item(102, 13).Value = -175.95
$ws.Cells.Item(102, 14).Value = -7794.4546
$ws.Cells.Item(132, 8).Value = 5452.086
$ws.Cells.Item(132, 9).Value = 5700.3794
$ws.Cells.Item(132, 11).Value = 17101.1382
$ws.Cells.Item(132, 13).Value = -14571.1382
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3450.5
$ws.Cells.Item(40, 9).Value = 3450.5
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 3450.5
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -3314.5
$ws.Cells.Item(100, 8).Value = 1945
$ws.Cells.Item(100, 9).Value = 1890
$ws.Cells.Item(100, 11).Value = 1890
$ws.Cells.Item(100, 13).Value = -1349
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 8070.6665
$ws.Cells.Item(74, 10).Value = 8070.6665
$ws.Cells.Item(74, 12).Value = 8070.6665
$ws.Cells.Item(74, 14).Value = -9942.666499999999
$ws.Cells.Item(77, 8).Value = 8070.6665
$ws.Cells.Item(77, 10).Value = 8070.6665
$ws.Cells.Item(77, 12).Value = 24211.9995
$ws.Cells.Item(77, 14).Value = -33571.99950000001
$ws.Cells.Item(100, 8).Value = 471
$ws.Cells.Item(100, 9).Value = 527.6667
$ws.Cells.Item(100, 10).Value = 301
$ws.Cells.Item(100, 11).Value = 1055.3334
$ws.Cells.Item(100, 12).Value = 602
$ws.Cells.Item(100, 13).Value = -514.3334
$ws.Cells.Item(100, 14).Value = -1684

# Cells whose trailing profit column no longer applies for this row
# (source row lost its HQ/secondary price entry) are cleared outright
# rather than zeroed, matching the upstream data for this refresh.
$wb.Worksheets.Item("ARM").Cells.Item(44, 14).ClearContents()
$wb.Worksheets.Item("LTW").Cells.Item(40, 14).ClearContents()

Write-Output "Applied Kujata profit-sheet refresh"
